# Add a new "CompactList" paragraph style, mirroring the existing
# "Compact" style: same base style (BodyText), same quick-style flag,
# and the same compact paragraph spacing (36 twips / 1.8pt before & after).
$d = $word.ActiveDocument

$compact = $d.Styles.Item("Compact")

$compactList = $d.Styles.Add("CompactList", 1)
$compactList.NameLocal = "Compact List"
$compactList.BaseStyle = $d.Styles.Item("BodyText")
$compactList.QuickStyle = $true
$compactList.ParagraphFormat.SpaceBefore = $compact.ParagraphFormat.SpaceBefore
$compactList.ParagraphFormat.SpaceAfter = $compact.ParagraphFormat.SpaceAfter
